# Update "想去人数" (number of people wanting to go) counts across sheets.
# Mirrors the gh-pages regenerated data snapshot (commit 456a3b4).

$wb = $excel.ActiveWorkbook

# -- Sheet "展览" (sheet1) --
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 37133
$ws1.Range("F4").Value = 631
$ws1.Range("F5").Value = 762
$ws1.Range("F12").Value = 521
$ws1.Range("F18").Value = 1153
$ws1.Range("F21").Value = 2479
$ws1.Range("F22").Value = 978
$ws1.Range("F23").Value = 546
$ws1.Range("F24").Value = 99
$ws1.Range("F25").Value = 1150
$ws1.Range("F29").Value = 1142

# -- Sheet "演出" (sheet2) --
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F4").Value = 327

# -- Sheet "全部类型" (sheet4) --
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value = 37133
$ws4.Range("F5").Value = 631
$ws4.Range("F6").Value = 762
$ws4.Range("F12").Value = 327
$ws4.Range("F17").Value = 521
$ws4.Range("F28").Value = 1153
$ws4.Range("F31").Value = 2479
$ws4.Range("F32").Value = 978
$ws4.Range("F33").Value = 546
$ws4.Range("F34").Value = 99
$ws4.Range("F35").Value = 1150
$ws4.Range("F40").Value = 1142
